$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows 375:385 (update fino a 20/09/2021)
$data = @(
    @(44449, 0, 0, 0),
    @(44450, 0, 0, 0),
    @(44451, 0, 0, 0),
    @(44452, 0, 0, 0),
    @(44453, 0, 0, 0),
    @(44454, 0, 0, 0),
    @(44455, 0, 0, 0),
    @(44456, 0, 0, 0),
    @(44457, 0, 0, 0),
    @(44458, 0, 0, 0),
    @(44459, 1, 1, 53.73455131649651)
)

$startRow = 375
$endRow = 385

# Copy formatting (style) of the last existing data row down into the new rows
$ws.Range("A374:D374").Copy()
$ws.Range("A375:D385").PasteSpecial(-4122) # xlPasteFormats

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
